$wb = $excel.ActiveWorkbook
$wsRef = $wb.Worksheets.Item("Reference")
$wsDef = $wb.Worksheets.Item("Defenses")

# Add the new "current match" rows (red_def2..5 / blue_def2..5) that were
# separated out from the Defenses reference table into their own rows.
$wsRef.Range("A27").Value = "red_def2"
$wsRef.Range("B27").Value = "red_def2"

$wsRef.Range("A28").Value = "red_def3"
$wsRef.Range("B28").Value = "red_def3"

$wsRef.Range("A29").Value = "red_def4"
$wsRef.Range("B29").Value = "red_def4"

$wsRef.Range("A30").Value = "red_def5"
$wsRef.Range("B30").Value = "red_def5"

$wsRef.Range("A31").Value = "blue_def2"
$wsRef.Range("B31").Value = "blue_def2"

$wsRef.Range("A32").Value = "blue_def3"
$wsRef.Range("B32").Value = "blue_def3"

$wsRef.Range("A33").Value = "blue_def4"
$wsRef.Range("B33").Value = "blue_def4"

$wsRef.Range("A34").Value = "blue_def5"
$wsRef.Range("B34").Value = "blue_def5"

# Update the selection on the Defenses sheet and move the active tab /
# selection over to the Reference sheet.
$wsDef.Range("B2:D9").Select()

$wsRef.Activate()
$wsRef.Range("B27").Select()
